# Updates the NATMI LR-pairs sheet with a fresh TPM re-run:
#  - rows 2-5 get recomputed values (same Sending/Ligand/Receptor/Target layout)
#  - rows 6-9 are newly appended (Lta cluster -> Tnfrsf14/ECs direction)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lta", "Tnfrsf14", "ECs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 4.901461666666667, 14.704385, 0.2124427850531459, 0.2124427850531459, 2.447562855276111, 22.028065697485, 0.2036234486942924, 0.2036234486942924),
    @("ECs", "Lta", "Tnfrsf14", "FAPs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 3.165799333333334, 9.497398, 0.1372144215401173, 0.1372144215401173, 1.580853505030889, 14.227681545278, 0.1315181107120274, 0.1315181107120274),
    @("ECs", "Lta", "Tnfrsf14", "MuSCs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 1.206743666666667, 3.620231, 0.05230357857032003, 0.05230357857032004, 0.6025918746767778, 5.423326872091001, 0.05013225111352747, 0.05013225111352748),
    @("ECs", "Lta", "Tnfrsf14", "Resolving-Mac", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 13.79790933333333, 41.393728, 0.5980392148364168, 0.5980392148364169, 6.890036617934221, 62.010329561408, 0.5732122526493622, 0.5732122526493623),
    @("Lta", "Tnfrsf14", "ECs", "ECs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 4.901461666666667, 14.704385, 0.2124427850531459, 0.2124427850531459, 0.1060088129266667, 0.95407931634, 0.00881933635885352, 0.008819336358853522),
    @("Lta", "Tnfrsf14", "ECs", "FAPs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 3.165799333333334, 9.497398, 0.1372144215401173, 0.1372144215401173, 0.06846990798133333, 0.616229171832, 0.0056963108280899, 0.005696310828089901),
    @("Lta", "Tnfrsf14", "ECs", "MuSCs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 1.206743666666667, 3.620231, 0.05230357857032003, 0.05230357857032004, 0.02609945202266666, 0.234895068204, 0.002171327456792558, 0.002171327456792558),
    @("Lta", "Tnfrsf14", "ECs", "Resolving-Mac", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 13.79790933333333, 41.393728, 0.5980392148364168, 0.5980392148364169, 0.2984211830613333, 2.685790647551999, 0.02482696218705461, 0.02482696218705461)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowValues[$c]
    }
}
